{"js": "// The report's paragraph describing \"\u05de\u05d8\u05dc\u05d4 0\" (\"assignment 0\") should read\n// \"...\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \u05d1\u05de\u05d8\u05dc\u05d4 0...\" instead of \"...\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \u05de\u05d8\u05dc\u05d4 0...\": insert the\n// single Hebrew letter \"\u05d1\" immediately before the word \"\u05de\u05d8\u05dc\u05d4\" in that\n// sentence. \"\u05de\u05d8\u05dc\u05d4\" also appears in the document title and later on, so\n// anchor the search on the unique (not-yet-edited) phrase \"\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \u05de\u05d8\u05dc\u05d4\"\n// to locate the exact word to fix, then insert \"\u05d1\" right before it.\nconst body = context.document.body;\n\nconst matches = body.search(\"\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \u05de\u05d8\u05dc\u05d4\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  const match = matches.items[0];\n  // Narrow the matched range down to just the \"\u05de\u05d8\u05dc\u05d4\" word at its end so we\n  // only insert \"\u05d1\" directly in front of it (leaving \"\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \" untouched).\n  const inner = match.search(\"\u05de\u05d8\u05dc\u05d4\", { matchCase: true, matchWholeWord: false });\n  inner.load(\"items\");\n  await context.sync();\n\n  const target = inner.items[0];\n  target.insertText(\"\u05d1\", \"Before\");\n  await context.sync();\n}\n", "ps1": "# The report's paragraph describing \"\u05de\u05d8\u05dc\u05d4 0\" (\"assignment 0\") should read\n# \"...\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \u05d1\u05de\u05d8\u05dc\u05d4 0...\" instead of \"...\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \u05de\u05d8\u05dc\u05d4 0...\": insert the\n# single Hebrew letter \"\u05d1\" immediately before the word \"\u05de\u05d8\u05dc\u05d4\" in that\n# sentence. \"\u05de\u05d8\u05dc\u05d4\" also appears elsewhere in the document (title, and later\n# \"\u05d4\u05d9\u05d5 \u05d1\u05de\u05d8\u05dc\u05d4 0\"), so anchor the Find on the unique phrase \"\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \u05de\u05d8\u05dc\u05d4\"\n# and replace it with \"\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \u05d1\u05de\u05d8\u05dc\u05d4\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \u05de\u05d8\u05dc\u05d4\"\n$find.Replacement.Text = \"\u05e9\u05d9\u05d9\u05e6\u05e8\u05e0\u05d5 \u05d1\u05de\u05d8\u05dc\u05d4\"\n$find.Forward = $true\n$find.Wrap = 0\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
